# Update countries & provincias Spain
#
# 1) Swap "Israel" / "Ucrania" rows (row 27 / row 28) so that the country
#    names line up with their (updated) statistics.
# 2) Update statistics for Ucrania (row 27), Israel (row 28),
#    Uzbekistan (row 62) and Taiwan (row 178).
# 3) Update the "Datos actualizados..." timestamp in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 27 / Row 28: swap country labels and refresh their figures -------
$ws.Range("A27").Value = "Ucrania"
$ws.Range("B27").Value = 309107
$ws.Range("C27").Value = 5469
$ws.Range("D27").Value = 129533
$ws.Range("E27").Value = 173788
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 113
$ws.Range("H27").Value = 5786

$ws.Range("A28").Value = "Israel"
$ws.Range("B28").Value = 304876
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = 278394
$ws.Range("E28").Value = 24219
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 2263

# --- Row 62: Uzbekistan -----------------------------------------------------
$ws.Range("B62").Value = 63737
$ws.Range("C62").Value = 214
$ws.Range("D62").Value = 60717
$ws.Range("E62").Value = 2487
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = 2
$ws.Range("H62").Value = 533

# --- Row 178: Taiwan ---------------------------------------------------------
$ws.Range("B178").Value = 543
$ws.Range("C178").Value = 3
$ws.Range("D178").Value = 495
$ws.Range("E178").Value = 41
$ws.Range("F178").Value = 0
$ws.Range("G178").Value = 0
$ws.Range("H178").Value = 7

# --- Updated timestamp -------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 20 de Octubre de 2020 a las 08:22"
